$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first two data rows (rows 2 and 3) were a "false start" - bogus test
# rows recorded before the real random-trading data collection began.
# Select them (as whole rows, the way a user would click-drag the row
# headers for rows 2 and 3) and delete them outright, which shifts the
# remaining rows (formerly 4 and 5) up to become the new rows 2 and 3.
$ws.Rows("2:3").Select()
$ws.Rows("2:3").Delete()
